$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the phone number column keeps being stored as text (not auto-converted to a number)
$ws.Range("H2:H7").NumberFormat = "@"

# Row 2 -> Mohammed Saneer (was row 4), A changes from 328 to 327
$ws.Range("A2").Value = 327
$ws.Range("B2").Value = "TE126"
$ws.Range("C2").Value = "Mohammed Saneer"
$ws.Range("D2").Value = "Mohammed"
$ws.Range("F2").Value = "Saneer"
$ws.Range("G2").Value = "saneer@dynasas.com"
$ws.Range("H2").Value = "971523235404"
$ws.Range("AE2").Value = "mohasane"
$ws.Range("AF2").Value = "c60019ad46409b66a812d6db61e2aa04f3ce644019d178761b5ef50e4af49db5"

# Row 3 -> Ajay Menon (was row 2), A changes from 329 to 328
$ws.Range("A3").Value = 328
$ws.Range("B3").Value = "TE125"
$ws.Range("C3").Value = "Ajay Menon"
$ws.Range("D3").Value = "Ajay"
$ws.Range("F3").Value = "Menon"
$ws.Range("G3").Value = "m.ajay@dynasas.com"
$ws.Range("H3").Value = "971508341694"
$ws.Range("AE3").Value = "ajaymeno"
$ws.Range("AF3").Value = "5aee76819f9b4633d11cd5abf7f8c2f6064ffa43e4cd72d1e091ef300d418008"

# Row 4 -> Ibrahim Rathwala (was row 3), A changes from 327 to 329
$ws.Range("A4").Value = 329
$ws.Range("B4").Value = "TE124"
$ws.Range("C4").Value = "Ibrahim Rathwala"
$ws.Range("D4").Value = "Ibrahim"
$ws.Range("F4").Value = "Rathwala"
$ws.Range("G4").Value = "yobozyt@gmail.com"
$ws.Range("H4").Value = "971561093935"
$ws.Range("AE4").Value = "ibrarath"
$ws.Range("AF4").Value = "30426da7ae9b01536a5b2b9ed461c17a323861583bd75284bc90a48065589658"

# Row 5 -> Ahmad Meda stays, only app username / password change
$ws.Range("AE5").Value = "ahmameda47"
$ws.Range("AF5").Value = "9f573e54f8f4ec488c5c8646a9cd42972847c16809e79e81aa13c2aa1064f956"

# Row 6 -> Hamza Ameen (was row 7), A changes from 325 to 326
$ws.Range("A6").Value = 326
$ws.Range("B6").Value = "TE127"
$ws.Range("C6").Value = "Hamza Ameen"
$ws.Range("D6").Value = "Hamza"
$ws.Range("F6").Value = "Ameen"
$ws.Range("G6").Value = "ameen.h@dynasas.com"
$ws.Range("H6").Value = "971566237633"
$ws.Range("AE6").Value = "hamzamee84"
$ws.Range("AF6").Value = "0dac17096b8e19a861230bdaafb6cf01a29c929c8a7b5ae61cb0a7c28275ed10"

# Row 7 -> Ruhban Gill (was row 6), A changes from 326 to 325
$ws.Range("A7").Value = 325
$ws.Range("B7").Value = "TE128"
$ws.Range("C7").Value = "Ruhban Gill"
$ws.Range("D7").Value = "Ruhban"
$ws.Range("F7").Value = "Gill"
$ws.Range("G7").Value = "Ruhban@dynasas.com"
$ws.Range("H7").Value = "971554980651"
$ws.Range("AE7").Value = "ruhbgill"
$ws.Range("AF7").Value = "34011906c331ed74ba1e3333a2a64a1bc4876c541a9c83c0b2e99637f386a236"
